$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  32"
$ws.Range("C9").Value = "Report Covering the Week  8/7/2023  Through  8/13/2023"

# --- Cells whose type/style changes (number<->text) ---
# Use NumberFormat="@" to force literal text storage where needed, then
# PasteSpecial(xlPasteFormats) from a sibling cell that already carries the
# exact target style, so the cells style index matches the target exactly.
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("C23").Value = 2
$ws.Range("C24").Copy() | Out-Null
$ws.Range("C23").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C26").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C27").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("D28").Value = 1
$ws.Range("D27").Copy() | Out-Null
$ws.Range("D28").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("E28").Value = -100
$ws.Range("E27").Copy() | Out-Null
$ws.Range("E28").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("D29").Value = 1
$ws.Range("D27").Copy() | Out-Null
$ws.Range("D29").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("E29").Value = -100
$ws.Range("E27").Copy() | Out-Null
$ws.Range("E29").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# --- Remaining numeric value updates (style unchanged) ---
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 9
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = 80
$ws.Range("L15").Value = 80
$ws.Range("M15").Value = 350
$ws.Range("N15").Value = -43.75
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -38.888888888888
$ws.Range("I16").Value = 96
$ws.Range("J16").Value = 122
$ws.Range("K16").Value = -21.311475409836
$ws.Range("L16").Value = 54.838709677419
$ws.Range("M16").Value = 2.127659574468
$ws.Range("N16").Value = -82.918149466192
$ws.Range("I17").Value = 141
$ws.Range("J17").Value = 163
$ws.Range("K17").Value = -13.496932515337
$ws.Range("L17").Value = 9.302325581395
$ws.Range("M17").Value = 123.809523809524
$ws.Range("N17").Value = -31.553398058252
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 127
$ws.Range("J18").Value = 145
$ws.Range("K18").Value = -12.413793103448
$ws.Range("L18").Value = 17.592592592592
$ws.Range("M18").Value = -16.993464052287
$ws.Range("N18").Value = -90.29793735676
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 28.571428571428
$ws.Range("F19").Value = 65
$ws.Range("G19").Value = 63
$ws.Range("H19").Value = 3.174603174603
$ws.Range("I19").Value = 420
$ws.Range("J19").Value = 474
$ws.Range("K19").Value = -11.392405063291
$ws.Range("L19").Value = 71.428571428571
$ws.Range("M19").Value = 83.406113537117
$ws.Range("N19").Value = -53.020134228187
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -62.5
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = -14.285714285714
$ws.Range("I20").Value = 120
$ws.Range("J20").Value = 118
$ws.Range("K20").Value = 1.694915254237
$ws.Range("L20").Value = 144.897959183673
$ws.Range("M20").Value = -2.439024390243
$ws.Range("N20").Value = -94.395142456795
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -10.81081081081
$ws.Range("F21").Value = 127
$ws.Range("G21").Value = 153
$ws.Range("H21").Value = -16.993464052287
$ws.Range("I21").Value = 917
$ws.Range("J21").Value = 1028
$ws.Range("K21").Value = -10.797665369649
$ws.Range("L21").Value = 53.344481605351
$ws.Range("M21").Value = 37.481259370314
$ws.Range("N21").Value = -82.121271202963
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 9
$ws.Range("K22").Value = 125
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = 12.5
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = -80
$ws.Range("I23").Value = 36
$ws.Range("J23").Value = 55
$ws.Range("K23").Value = -34.545454545454
$ws.Range("L23").Value = -10
$ws.Range("M23").Value = 111.764705882353
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -25.806451612903
$ws.Range("F24").Value = 81
$ws.Range("G24").Value = 158
$ws.Range("H24").Value = -48.734177215189
$ws.Range("I24").Value = 755
$ws.Range("J24").Value = 954
$ws.Range("K24").Value = -20.859538784067
$ws.Range("L24").Value = 21.774193548387
$ws.Range("M24").Value = 18.897637795275
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = -38.461538461538
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 43
$ws.Range("H25").Value = -30.232558139534
$ws.Range("I25").Value = 250
$ws.Range("J25").Value = 289
$ws.Range("K25").Value = -13.494809688581
$ws.Range("L25").Value = 20.772946859903
$ws.Range("M25").Value = -22.360248447205
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 11
$ws.Range("K26").Value = 9.090909090909
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = -70
$ws.Range("J27").Value = 27
$ws.Range("K27").Value = 3.703703703703
$ws.Range("L27").Value = 16.666666666666
$ws.Range("G28").Value = 3
$ws.Range("J28").Value = 8
$ws.Range("K28").Value = -62.5
$ws.Range("G29").Value = 2
$ws.Range("J29").Value = 7
$ws.Range("K29").Value = -57.142857142857
$ws.Range("F30").Value = 2
$ws.Range("H30").Value = 100
$ws.Range("I30").Value = 10
$ws.Range("K30").Value = -37.5
$ws.Range("L30").Value = 66.666666666666
